$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New prediction data for the 2021-01-09 forecast, re-run with weather data,
# appended as rows 68-75 (weeks 10 Jan 2021 .. 06 Mar 2021)

$dayMade = "2021-01-09"
$model = "KNN"

$weeks = @(
    "10 Jan -- 16 Jan 2021",
    "17 Jan -- 23 Jan 2021",
    "24 Jan -- 30 Jan 2021",
    "31 Jan -- 06 Feb 2021",
    "07 Feb -- 13 Feb 2021",
    "14 Feb -- 20 Feb 2021",
    "21 Feb -- 27 Feb 2021",
    "28 Feb -- 06 Mar 2021"
)

$predictions = @(860.29, 862.0700000000001, 826.12, 794.46, 766.6900000000001, 767.8099999999999, 802.62, 904.74)

$startRow = 68

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $r = $startRow + $i
    # Column A holds a literal text string ("2021-01-09"), not a real date,
    # matching the rest of the sheet -- force text so Excel doesn't
    # auto-convert it to a date serial, then restore the default style.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $dayMade
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $weeks[$i]
    $ws.Cells.Item($r, 4).Value = $predictions[$i]
    $ws.Cells.Item($r, 6).Value = $model
}

# First new row also carries Real/difference/MAE-related figures
$ws.Cells.Item(68, 3).Value = 3333.57
$ws.Cells.Item(68, 5).Value = 2473.28
$ws.Cells.Item(68, 9).Value = 2375.32
$ws.Cells.Item(68, 10).Value = 2150.77
$ws.Cells.Item(68, 11).Value = 69.94
